$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "03828edbbdf81522657d2618d06ddf1bd611fa0abd1d4cc230fc6c17b48a0b3709"
$ws.Range("E3").Value = "03ab5f057da235b71b082f906292a37aa35e4d253b53618a7ebd9c2cf059f54e75"
$ws.Range("E4").Value = "02b4adb050cc2429a8375c8a5d0a133952d8edc20900c7d4dd8f5ad1cf9113dd58"
$ws.Range("E5").Value = "021d85da863f80e8067960374d333aaaad9f0cffb51f453fa2d20cc1b49bbbd8b3"
$ws.Range("E6").Value = "02502dbc748bded2563cf985b6560073c2fee73251950ec58534194dfcdbe1662e"

$ws.Range("E6").Select()
